$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = 9440.75
$ws.Range("B17").Value = 9798.39
$ws.Range("C17").Value = 277
$ws.Range("D17").Value = 287.11
$ws.Range("E17").Value = $true
$ws.Range("F17").Value = 3.65
$ws.Range("G17").Value = 42626.545347222222
$ws.Range("G17").NumberFormat = "m/d/yy h:mm"
$ws.Range("H17").Value = $false
